$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the raw observation counts used by the logistic fit ---
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 13

# --- Update the solved (Solver-fitted) mu/s parameters for the logistic curve ---
# (these are the converged values produced by re-running Solver on D5 = SUMSQ(D2:D4)
#  against $G$1:$G$2 with the updated data above)
$ws.Range("G1").Value = 18.007274784929372
$ws.Range("G2").Value = 4.5890312131293562

# --- Remove the leftover scratch/backup copies of the solver output that were
#     parked off to the side (rows 9, 10, 13, 14, and the "frombelow"/"fromabove"
#     labels on rows 13/17) ---
$ws.Range("H9").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("G13").ClearContents()
$ws.Range("H13").ClearContents()
$ws.Range("H14").ClearContents()
$ws.Range("G17").ClearContents()
$ws.Range("H17").ClearContents()
$ws.Range("H18").ClearContents()

# --- Update the Solver add-in's saved engine/version settings for this run ---
$estName = $ws.Names.Add("solver_est", "=1")
$estName.Visible = $false
$nwtName = $ws.Names.Add("solver_nwt", "=1")
$nwtName.Visible = $false

foreach ($n in $ws.Names) {
  if ($n.Name() -eq "Sheet1!solver_ver") {
    $n.RefersTo = "=3"
  }
}

# --- Move the active selection to where the user last left off ---
$null = $ws.Range("G6").Select()

Write-Output "edit applied"
